$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update depth From/To values (columns B and C) ---
$ws.Range("C11").Value = 6.75
$ws.Range("C12").Value = 7.65
$ws.Range("B14").Value = 9.3
$ws.Range("C14").Value = 9.35
$ws.Range("B15").Value = 9.35
$ws.Range("C15").Value = 10.5
$ws.Range("B16").Value = 10.5
$ws.Range("C17").Value = 12.3
$ws.Range("C18").Value = 13.75
$ws.Range("C20").Value = 16.55
$ws.Range("C21").Value = 17.1
$ws.Range("B22").Value = 17.1
$ws.Range("C23").Value = 18.3
$ws.Range("B24").Value = 18.35

# --- Update N column (blow count) values ---
$ws.Range("N9").Value = 39
$ws.Range("N10").Value = 40
$ws.Range("N11").Value = 240
$ws.Range("N12").Value = 240
$ws.Range("N13").Value = 300
$ws.Range("N14").Value = 2000
$ws.Range("N15").Value = 2000
$ws.Range("N16").Value = 2000
$ws.Range("N17").Value = 120
$ws.Range("N18").Value = 240
$ws.Range("N19").Value = 200
$ws.Range("N20").Value = 1200
$ws.Range("N21").Value = 2000
$ws.Range("N22").Value = 2000
$ws.Range("N23").Value = 2000
$ws.Range("N24").Value = 2000
$ws.Range("N25").Value = 2000

# --- Update selected cell / active selection ---
$ws.Range("N17").Select() | Out-Null
